$d = $word.ActiveDocument

function New-XmlPackage([string]$bodyFragment) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $bodyFragment + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# --- Change 1 --------------------------------------------------------
# Collapse the run-fragmented "employer would then follow a link..."
# paragraph down to two runs (one red "The employer would then follow a
# link " run, one plain "to the login page..." run), matching the other
# (graduate) paragraph's structure.
$employerLinkPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*would then follow a link*Employer dashboard*") {
        $employerLinkPara = $p
    }
}

$body1 = '<w:body><w:p><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">The employer would then follow a link </w:t></w:r><w:r><w:t>to the login page, login using the given username and password and be directed the Employer dashboard.</w:t></w:r></w:p></w:body>'
$employerLinkRange = $d.Range($employerLinkPara.Range.Start, $employerLinkPara.Range.End - 1)
$employerLinkRange.InsertXML((New-XmlPackage $body1))

# --- Change 2 --------------------------------------------------------
# a) Strip the red-colour formatting from the "They could see who is
#    graduate, employer or TDA using a flag." paragraph (both the
#    paragraph-mark run properties and the text run's properties).
$flagPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*They could see who is graduate*") {
        $flagPara = $p
    }
}

$body2 = '<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>They could see who is graduate, employer or TDA using a flag.</w:t></w:r></w:p></w:body>'
$flagPara.Range.InsertXML((New-XmlPackage $body2))

# reacquire the paragraph after replacing its XML
$flagPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*They could see who is graduate*") {
        $flagPara = $p
    }
}

# b) Insert a brand-new paragraph straight after it describing the new
#    "Register new graduate, employer or tda admin" user story, with the
#    same list formatting/red colour as the other TDA admin bullets.
$flagPara.Range.InsertParagraphAfter()

$newPara = $null
$after = $false
foreach ($p in $d.Paragraphs) {
    if ($after) { $newPara = $p; break }
    if ($p.Range.Text -like "*They could see who is graduate*") {
        $after = $true
    }
}

$body3 = '<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">Register new graduate, employer or </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>tda</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve"> admin &#8211; username, password</w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>roles</w:t></w:r></w:p></w:body>'
$newPara.Range.InsertXML((New-XmlPackage $body3))
